# Update countries & provincias Spain
# Refresh the COVID country table to the 13:31 snapshot and re-sort by
# total cases (column B) descending, which is how "Kuwait" overtakes
# "Emiratos Arabes Unidos" in the ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 13:31"

# --- 2. Refresh per-country figures (Casos totales, Nuevos casos, Casos
#         activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 3220559
$ws.Range("C4").Value = 560
$ws.Range("E4").Value = 1658248

# Iran (row 13)
$ws.Range("B13").Value = 252720
$ws.Range("C13").Value = 2262
$ws.Range("D13").Value = 215015
$ws.Range("E13").Value = 25258
$ws.Range("G13").Value = 142
$ws.Range("H13").Value = 12447

# Alemania (row 19)
$ws.Range("B19").Value = 199216
$ws.Range("C19").Value = 18
$ws.Range("E19").Value = 6491

# Bielorrusia (row 31)
$ws.Range("B31").Value = 64604
$ws.Range("C31").Value = 193
$ws.Range("D31").Value = 54254
$ws.Range("E31").Value = 9896
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = 454

# Kuwait (row 38, before re-sort)
$ws.Range("B38").Value = 53580
$ws.Range("C38").Value = 740
$ws.Range("D38").Value = 43214
$ws.Range("E38").Value = 9983
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 383

# Suiza (row 49)
$ws.Range("B49").Value = 32690
$ws.Range("C49").Value = 104
$ws.Range("E49").Value = 1324

# Nepal (row 64)
$ws.Range("B64").Value = 16649
$ws.Range("C64").Value = 118
$ws.Range("D64").Value = 8011
$ws.Range("E64").Value = 8603

# Australia (row 73)
$ws.Range("B73").Value = 9359
$ws.Range("C73").Value = 300
$ws.Range("E73").Value = 1677

# Estado de Palestina (row 93)
$ws.Range("B93").Value = 5551
$ws.Range("C93").Value = 331
$ws.Range("D93").Value = 536
$ws.Range("E93").Value = 4988

# Madagascar (row 99)
$ws.Range("B99").Value = 4143
$ws.Range("C99").Value = 361
$ws.Range("D99").Value = 2183
$ws.Range("E99").Value = 1926
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 34

# Malta (row 154)
$ws.Range("D154").Value = 658
$ws.Range("E154").Value = 7

# Siria (row 160)
$ws.Range("B160").Value = 394
$ws.Range("C160").Value = 22
$ws.Range("E160").Value = 252
$ws.Range("G160").Value = 2
$ws.Range("H160").Value = 16

# --- 3. Re-sort the country table by total cases (column B) descending,
#         now that Kuwait's updated total (53580) passes Emiratos Arabes
#         Unidos (53577) ---
$sortRange = $ws.Range("A4:H219")
$keyRange = $ws.Range("B4:B219")
$sortRange.Sort($keyRange, 2, $null, $null, 1, $null, 1, 0)
